# daily auto push: 2026-01-09 02:27 UTC
#
# The sheet is a log of date/day-of-week/hour/ranking rows. A new
# measurement for 2026/01/09 (already the date of the last existing row,
# row 584) needs to be appended right after the existing 2026/01/09 row,
# pushing every following row down by one. The new row re-uses the same
# date ("2026/01/09") and weekday ("金") as row 584, with new C/D values
# of 8 and 23.
#
# To make sure the new cell A585 stays a plain text value "2026/01/09"
# (and not get auto-converted by Excel into a date serial number, which
# happens if you just assign the string to .Value), we copy the existing
# row 584 (which is already correctly typed) down into a freshly
# inserted row 585, and then only overwrite the C/D (time/ranking)
# values that differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "2026/01/09" row and insert it as a new row at 585,
# shifting row 585 (and everything after it) down by one row.
$ws.Range("A584:D584").Copy()
$ws.Range("A585").Insert(-4121)   # -4121 = xlShiftDown

# Update the time / ranking columns of the newly inserted row.
$ws.Range("C585").Value = 8
$ws.Range("D585").Value = 23
